$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old stray cells that are being removed (C18, L25, L26)
$ws.Range("C18").ClearContents()
$ws.Range("L25").ClearContents()
$ws.Range("L26").ClearContents()

# Fill in new column D (sex = "U") and E (environ = 1) for rows 2-7
# plus F/G (sire/dam) values that are new/changed for rows 4-7
$ws.Range("D2").Value = "U"
$ws.Range("E2").Value = 1

$ws.Range("D3").Value = "U"
$ws.Range("E3").Value = 1

$ws.Range("D4").Value = "U"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 2
$ws.Range("G4").Value = 1

$ws.Range("D5").Value = "U"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 3
$ws.Range("G5").Value = 3

$ws.Range("D6").Value = "U"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 3
$ws.Range("G6").Value = 1

$ws.Range("D7").Value = "U"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 3
$ws.Range("G7").Value = 2

# Update sheet view selection
$ws.Range("K11").Select()
